$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individual Costs")

# row, Vehicle, Electricity_Cost, Degradation_Cost, GHG_Cost, Total Charge
$data = @(
    @(2, "P_1092", 2705.384570327582, 156.3399505236742, 91.07593037516521, 7090.247189282279),
    @(3, "P_1307", 1326.101245, 87.30917999999996, 53.283533436125, 3959.599999999998),
    @(4, "P_1111", 1184.3886825, 78.49733850000007, 48.64042827868546, 3559.970000000003),
    @(5, "P_1393", 1861.425984098608, 122.03978091017, 52.60679542982393, 5534.68394150431),
    @(6, "P_1304", 955.8107446153845, 61.56360000000002, 36.57289490219227, 2792),
    @(7, "P_1279", 2936.149079999988, 168.40908, 96.05876624089991, 7637.599999999999),
    @(8, "P_1371", 1180.515854016911, 56.64827908777654, 33.13819390409638, 2569.0829518266),
    @(9, "P_1419", 1232.568295271651, 67.53164665920092, 38.6483330356762, 3062.659712435416),
    @(10, "P_1421", 1157.107673240658, 65.61978126510103, 36.30818143506726, 2975.953798870795),
    @(11, "P_1100", 2354.434530751681, 149.7956827499998, 89.23199304412083, 6793.454999999992),
    @(12, "P_1141", 3568.465971981511, 178.5829499999998, 102.1982032613868, 8098.999999999988),
    @(13, "P_1093", 1903.013565189496, 123.3650092500005, 76.84668256242304, 5594.785000000023),
    @(14, "P_1257", 1575.774973634665, 88.96876707959395, 53.72884433391435, 4034.864720162991),
    @(15, "P_1143", 1760.736591886139, 100.4333399999999, 59.91512060420938, 4554.799999999996),
    @(16, "P_1272", 2925.324735442888, 156.9272040000004, 94.16916893510829, 7116.880000000015),
    @(17, "P_1127", 3447.003767666663, 204.9613649999997, 123.2169825691001, 9295.299999999983),
    @(18, "P_1260", 2558.516710438244, 145.7358503319236, 79.99428898205711, 6609.335615960253),
    @(19, "P_1123", 3313.07877894633, 138.0354095644854, 74.6646107620076, 6260.109277300924),
    @(20, "P_1295", 1347.933576000001, 88.51399199999986, 54.93293609434996, 4014.239999999993),
    @(21, "P_1375", 2300.259317575762, 120.1460400000001, 61.20184795059098, 5448.800000000005),
    @(22, "P_1376", 1549.559821681081, 93.68710783880661, 47.30750649483431, 4248.848428063791),
    @(23, "P_1414", 762.864850128035, 45.86764546122762, 28.6187456631727, 2080.165327039801),
    @(24, "P_1131", 2088.648187116987, 108.7814700000004, 61.0705177994389, 4933.400000000017),
    @(25, "P_1112", 2262.991851750003, 96.61075200000002, 44.59591799591247, 4381.440000000001),
    @(26, "P_1132", 1257.320938702072, 70.53331950000009, 41.72447031467404, 3198.790000000004),
    @(27, "P_1253", 1078.5775875, 54.93316500000002, 35.12202847676254, 2491.300000000001),
    @(28, "P_1135", 1868.35710336842, 118.6907400000003, 71.40030787805004, 5382.800000000015),
    @(29, "P_1109", 2636.903667083329, 158.2335562500003, 98.31126052625729, 7176.125000000015),
    @(30, "P_1353", 1187.458966537751, 64.09599683416663, 40.58505108367701, 2906.847928987148),
    @(31, "P_1424", 367.5125547876154, 19.65511432815084, 12.26175685510235, 891.3884049048),
    @(32, "P_1427", 1648.670986013978, 98.43651703802752, 58.59265312081038, 4464.241135511452),
    @(33, "P_1217", 601.1749396710526, 25.474365, 11.32326814729268, 1155.3),
    @(34, "P_1271", 3616.056223939393, 190.32678, 100.2043873143121, 8631.599999999999),
    @(35, "P_1368", 1383.290288718546, 74.43648098050134, 46.05181330529004, 3375.804126099834),
    @(36, "P_1094", 2328.450095945075, 114.314256, 66.33131540430304, 5184.319999999997),
    @(37, "P_1285", 1429.641299999999, 86.91228000000004, 50.52206982505907, 3941.600000000001),
    @(38, "P_1288", 650.2027944184783, 38.28530474999997, 21.60099836996058, 1736.294999999999),
    @(39, "P_1357", 1638.316774720954, 102.2198741718037, 65.1908425176393, 4635.821957904929),
    @(40, "P_1280", 1096.869247090908, 66.38461199999999, 39.61473453095603, 3010.639999999999),
    @(41, "P_1281", 1932.171259423079, 117.436095, 68.50094724358814, 5325.900000000001),
    @(42, "P_1296", 1751.640419831837, 88.03021499999994, 52.77573966803654, 3992.299999999997),
    @(43, "P_1367", 1422.875604374262, 83.50570956409045, 50.09945291411204, 3787.107009709317),
    @(44, "P_1125a", 970.9552616568085, 41.83201878493561, 20.95130077938576, 1897.143709067375),
    @(45, "P_1087", 3167.998191812507, 181.5042442499993, 108.871073795229, 8231.484999999966),
    @(46, "P_1422", 847.5148239900864, 49.29635530629912, 28.5739729539983, 2235.66237216776),
    @(47, "P_1137", 1611.088021569844, 78.35881930843867, 41.2342618612266, 3553.687950496084),
    @(48, "P_1294", 735.9627238962704, 35.58278763351645, 18.68265108277637, 1613.731865465599),
    @(49, "P_1091", 5574.829652889717, 325.5459794999994, 201.0975713303817, 14763.98999999997),
    @(50, "P_1125", 878.078724423309, 37.52948423153453, 21.83733499788648, 1702.017425466418),
    @(51, "P_1098", 2294.152660798979, 121.5743120129884, 58.89065695678865, 5513.574240951853)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
    $ws.Cells.Item($r, 6).Value = $item[5]
}
